# ============================================================
# MAJ automatique BRVM - mise a jour des feuilles Recommandations
# et Top_YTD suite au recalcul quotidien des indices/valeurs.
# ============================================================

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# ---- Feuille 'Recommandations' : lignes 2 a 40 (tri par Variation Totale desc.) ----
# Ligne 2
$ws1.Range("D2").Value = 677.92
$ws1.Range("E2").Value = 168.79
# Ligne 3
$ws1.Range("D3").Value = 601.01
$ws1.Range("E3").Value = 150.88
# Ligne 4
$ws1.Range("D4").Value = 583.22
$ws1.Range("E4").Value = 145.78
# Ligne 5
$ws1.Range("D5").Value = 581.53
$ws1.Range("E5").Value = 145.52
# Ligne 6
$ws1.Range("D6").Value = 466.42
$ws1.Range("E6").Value = 118.54
# Ligne 7
$ws1.Range("A7").Value = "BRVM - ENERGIE"
$ws1.Range("C7").Value = 4
$ws1.Range("D7").Value = 462.68
$ws1.Range("E7").Value = 117.38
# Ligne 8
$ws1.Range("A8").Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Range("D8").Value = 376.99
$ws1.Range("E8").Value = 94.26
# Ligne 9
$ws1.Range("A9").Value = "BRVM - CONSOMMATION DE BASE     (**)"
$ws1.Range("C9").Value = 1
$ws1.Range("D9").Value = 231.36
$ws1.Range("E9").Value = 231.36
# Ligne 10
$ws1.Range("A10").Value = "BRVM-PRINCIPAL     (**)"
$ws1.Range("C10").Value = 1
$ws1.Range("D10").Value = 228.7
$ws1.Range("E10").Value = 228.7
# Ligne 11
$ws1.Range("A11").Value = "BRVM-PRINCIPAL"
$ws1.Range("C11").Value = 1
$ws1.Range("D11").Value = 227.6
$ws1.Range("E11").Value = 227.6
# Ligne 12
$ws1.Range("A12").Value = "BRVM - CONSOMMATION DE BASE"
$ws1.Range("D12").Value = 225.22
$ws1.Range("E12").Value = 225.22
# Ligne 13
$ws1.Range("A13").Value = "BRVM – COMPOSITE TOTAL RETURN     (**)"
$ws1.Range("D13").Value = 136
$ws1.Range("E13").Value = 136
# Ligne 15
$ws1.Range("B15").Value = 4
$ws1.Range("D15").Value = 28.94
$ws1.Range("E15").Value = 6.89
# Ligne 17
$ws1.Range("A17").Value = "SAFCA CI (SAFC)"
$ws1.Range("C17").Value = 0
$ws1.Range("D17").Value = 11.31
$ws1.Range("E17").Value = 5.57
$ws1.Range("G17").Value = "➖ Neutre"
# Ligne 18
$ws1.Range("A18").Value = "SUCRIVOIRE (SCRC)"
$ws1.Range("B18").Value = 1
$ws1.Range("C18").Value = 0
$ws1.Range("D18").Value = 7.27
$ws1.Range("E18").Value = 7.27
$ws1.Range("G18").Value = "➖ Neutre"
# Ligne 19
$ws1.Range("A19").Value = "CORIS BANK INTERNATIONAL (CBIBF)"
$ws1.Range("D19").Value = 7.11
$ws1.Range("E19").Value = 7.11
# Ligne 21
$ws1.Range("A21").Value = "VIVO ENERGY CI (SHEC)"
$ws1.Range("D21").Value = 4.46
$ws1.Range("E21").Value = 4.46
# Ligne 23
$ws1.Range("A23").Value = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$ws1.Range("D23").Value = 3.35
$ws1.Range("E23").Value = 3.35
# Ligne 24
$ws1.Range("A24").Value = "ONATEL BF (ONTBF)"
$ws1.Range("C24").Value = 1
$ws1.Range("D24").Value = 3.23
$ws1.Range("E24").Value = -1.15
$ws1.Range("G24").Value = "👀 À surveiller"
# Ligne 25
$ws1.Range("A25").Value = "SMB CI (SMBC)"
$ws1.Range("C25").Value = 0
$ws1.Range("D25").Value = 2.99
$ws1.Range("E25").Value = 2.99
$ws1.Range("G25").Value = "➖ Neutre"
# Ligne 26
$ws1.Range("A26").Value = "UNILEVER CI (UNLC)"
$ws1.Range("C26").Value = 1
$ws1.Range("D26").Value = 1.24
$ws1.Range("E26").Value = 7.49
$ws1.Range("G26").Value = "👀 À surveiller"
# Ligne 27
$ws1.Range("A27").Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Range("B27").Value = 2
$ws1.Range("C27").Value = 2
$ws1.Range("D27").Value = 0.4
$ws1.Range("E27").Value = 4.55
# Ligne 28
$ws1.Range("A28").Value = "EVIOSYS PACKAGING SIEM CI (SEMC)"
$ws1.Range("D28").Value = -0.32
$ws1.Range("E28").Value = -6.81
# Ligne 30
$ws1.Range("A30").Value = "NESTLE CI (NTLC)"
$ws1.Range("D30").Value = -0.89
$ws1.Range("E30").Value = -0.89
# Ligne 31
$ws1.Range("A31").Value = "SAPH CI (SPHC)"
$ws1.Range("D31").Value = -0.98
$ws1.Range("E31").Value = -0.98
# Ligne 32
$ws1.Range("A32").Value = "BANK OF AFRICA NG (BOAN)"
$ws1.Range("D32").Value = -1.14
$ws1.Range("E32").Value = -1.14
# Ligne 33
$ws1.Range("A33").Value = "SOGB CI (SOGC)"
$ws1.Range("D33").Value = -2.04
$ws1.Range("E33").Value = -2.04
# Ligne 34
$ws1.Range("A34").Value = "ORAGROUP TOGO (ORGT)"
$ws1.Range("D34").Value = -2.08
$ws1.Range("E34").Value = -2.08
# Ligne 35
$ws1.Range("A35").Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$ws1.Range("B35").Value = 0
$ws1.Range("C35").Value = 1
$ws1.Range("D35").Value = -2.89
$ws1.Range("E35").Value = -2.89
$ws1.Range("G35").Value = "➖ Neutre"
# Ligne 36
$ws1.Range("C36").Value = 1
$ws1.Range("D36").Value = -3.67
# Ligne 37
$ws1.Range("A37").Value = "SETAO CI (STAC)"
$ws1.Range("C37").Value = 1
$ws1.Range("D37").Value = -4
$ws1.Range("E37").Value = -4
# Ligne 38
$ws1.Range("A38").Value = "CFAO MOTORS CI (CFAC)"
$ws1.Range("C38").Value = 3
$ws1.Range("D38").Value = -4.9
$ws1.Range("E38").Value = -1.43
# Ligne 39
$ws1.Range("A39").Value = "SOLIBRA CI (SLBC)"
$ws1.Range("B39").Value = 0
$ws1.Range("C39").Value = 1
$ws1.Range("D39").Value = -5.22
$ws1.Range("E39").Value = -5.22
$ws1.Range("F39").Value = "🟡 Observer"
$ws1.Range("G39").Value = "➖ Neutre"
# Ligne 40
$ws1.Range("A40").Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws1.Range("B40").Value = 0
$ws1.Range("C40").Value = 2
$ws1.Range("D40").Value = -5.44
$ws1.Range("E40").Value = -2.63
$ws1.Range("F40").Value = "🟡 Observer"
$ws1.Range("G40").Value = "➖ Neutre"

# ---- Feuille 'Top_YTD' : lignes 2 a 11 (colonne B = Progression YTD) ----
# Ligne 2
$ws2.Range("B2").Value = 5173.3
# Ligne 3
$ws2.Range("B3").Value = 3822.03
# Ligne 4
$ws2.Range("B4").Value = 3550.58
# Ligne 5
$ws2.Range("B5").Value = 3525.54
# Ligne 6
$ws2.Range("B6").Value = 2101.06
# Ligne 7
$ws2.Range("B7").Value = 2063.39
# Ligne 8
$ws2.Range("B8").Value = 1323.7
# Ligne 9
$ws2.Range("B9").Value = 231.36
# Ligne 10
$ws2.Range("B10").Value = 228.7
# Ligne 11
$ws2.Range("A11").Value = "BRVM-PRINCIPAL"
$ws2.Range("B11").Value = 227.6
